# Automatische test-sync: 2025-08-05 18:37:50
# Append a new log entry (row 36) to the "Logs" sheet and bump the
# matching "Planning / Afspraak" tally on the "Dashboard" sheet.

$wb = $excel.ActiveWorkbook
$logs = $wb.Worksheets.Item("Logs")
$dashboard = $wb.Worksheets.Item("Dashboard")

$newRow = 36

$logs.Cells.Item($newRow, 1).Value = "Leg dit even neer bij Koen."
$logs.Cells.Item($newRow, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item($newRow, 3).Value = "Testmail #15: Leg dit even neer bij Koen."
$logs.Cells.Item($newRow, 4).Value = "Planning / Afspraak"
$logs.Cells.Item($newRow, 5).Value = "Bedankt, we hebben dit doorgestuurd naar planning@bedrijf.nl."
$logs.Cells.Item($newRow, 6).Value = "2025-08-05 18:37:11"
$logs.Cells.Item($newRow, 7).Value = "Ja"
$logs.Cells.Item($newRow, 8).Value = "Ja"
$logs.Cells.Item($newRow, 9).Value = "Nee"
$logs.Cells.Item($newRow, 10).Value = "Nee"

# Extend the conditional-formatting ranges (D/G/H/I/J 2:35 -> 2:36) so the
# newly-added row is covered just like the rest of the log table. Modifying
# one rule's applies-to range re-seats the whole sibling rule group sharing
# that sqref.
$logs.Range("D2:D35").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D36"))
$logs.Range("G2:G35").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G36"))
$logs.Range("H2:H35").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H36"))
$logs.Range("I2:I35").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I36"))
$logs.Range("J2:J35").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("J2:J36"))

# The new row's category is "Planning / Afspraak", so bump its Dashboard count.
$dashboard.Range("B2").Value = 19
